$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 182 - this pushes the existing rows 182..217 down to 183..218,
# carrying their formatting (including the date style on column D) along with them.
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row 182 with the new weekly record.
$ws.Cells.Item(182, 1).Value = 8
$ws.Cells.Item(182, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(182, 3).Value = "Coquimbo"
$ws.Cells.Item(182, 4).Value = 45015
$ws.Cells.Item(182, 5).Value = 4
$ws.Cells.Item(182, 6).Value = 100112040
$ws.Cells.Item(182, 7).Value = "Cilantro"
$ws.Cells.Item(182, 8).Value = "Sin especificar"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 2400
$ws.Cells.Item(182, 11).Value = 1800
$ws.Cells.Item(182, 12).Value = 2000
$ws.Cells.Item(182, 13).Value = 1900
$ws.Cells.Item(182, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(182, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(182, 16).Value = 1267
$ws.Cells.Item(182, 17).Value = 1.5
$ws.Cells.Item(182, 18).Value = "Hortaliza"
